# [DEV] Initial steps CMMI
# Populates Plan3 ("Plan3" / sheet3) with the CMMI process-area level table,
# formats it, merges the level cells, and makes Plan3 the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- constants (Excel enum values) ----
$xlCenter  = -4108
$xlRight   = -4152
$xlJustify = -4130

# ---- header row (row 1) ----
$ws.Cells.Item(1, 1).Value = "Nível"
$ws.Cells.Item(1, 2).Value = "Sigla"
$ws.Cells.Item(1, 3).Value = "Descrição"

$hdr = $ws.Range("A1:C1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 12
$hdr.Font.Bold = $true
$hdr.VerticalAlignment = $xlCenter

$ws.Range("A1").HorizontalAlignment = $xlCenter
$ws.Range("B1").HorizontalAlignment = $xlRight
$ws.Range("B1").IndentLevel = 1
$ws.Range("C1").HorizontalAlignment = $xlJustify

$ws.Rows.Item(1).RowHeight = 15.75

# ---- data blocks ----
# level 2 -> rows 2-8
$block1 = @(
    @("CM", "- Gestão da configuração"),
    @("MA", "- Medição e Análise"),
    @("PMC", "- Projeto de Monitoramento e Controle"),
    @("PP", "- Planejamento de Projetos"),
    @("PPQA", "- Processo e Produto Quality Assurance"),
    @("REQM", "- Gerenciamento de Requisitos"),
    @("SAM", "- Gerenciamento de acordo com o fornecedor")
)

# level 3 -> rows 10-20
$block2 = @(
    @("DAR", "- Análise de decisão e resolução"),
    @("IPM", "- Gestão Integrada de Projetos"),
    @("OPD", "- Definição do Processo Organizacional"),
    @("OPF", "- Foco no Processo Organizacional "),
    @("OT", "- Treinamento Organizacional"),
    @("PI", "- Integração de Produto"),
    @("RD", "- Desenvolvimento de Requisitos"),
    @("RSKM", "- Gestão de Riscos"),
    @("TS", "- Solução Técnica"),
    @("VAL", "- Validação"),
    @("VER", "- Verificação")
)

# level 4 -> rows 22-23
$block3 = @(
    @("OPP ", "- Performance do Processo Organizacional"),
    @("QPM ", "- Projeto quantativamente gerenciado")
)

# level 5 -> rows 25-26
$block4 = @(
    @("CAR ", "- Análise e Resolução de Causas"),
    @("OPM ", "- Gestão de Desempenho Organizacional")
)

# ---- write block 1 (rows 2-8), level number 2 in A2 ----
$ws.Cells.Item(2, 1).Value = 2
for ($i = 0; $i -lt $block1.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 2).Value = $block1[$i][0]
    $ws.Cells.Item($r, 3).Value = $block1[$i][1]
}

# ---- blank spacer row 9 ----
$ws.Cells.Item(9, 1).Value = ""

# ---- write block 2 (rows 10-20), level number 3 in A10 ----
$ws.Cells.Item(10, 1).Value = 3
for ($i = 0; $i -lt $block2.Length; $i++) {
    $r = 10 + $i
    $ws.Cells.Item($r, 2).Value = $block2[$i][0]
    $ws.Cells.Item($r, 3).Value = $block2[$i][1]
}

# ---- blank spacer row 21 ----
$ws.Cells.Item(21, 1).Value = ""

# ---- write block 3 (rows 22-23), level number 4 in A22 ----
$ws.Cells.Item(22, 1).Value = 4
for ($i = 0; $i -lt $block3.Length; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 2).Value = $block3[$i][0]
    $ws.Cells.Item($r, 3).Value = $block3[$i][1]
}

# ---- blank spacer row 24 ----
$ws.Cells.Item(24, 1).Value = ""

# ---- write block 4 (rows 25-26), level number 5 in A25 ----
$ws.Cells.Item(25, 1).Value = 5
for ($i = 0; $i -lt $block4.Length; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 2).Value = $block4[$i][0]
    $ws.Cells.Item($r, 3).Value = $block4[$i][1]
}

# ---- formatting: column A (level numbers), rows 2-8 and 10-20 ----
$colA_2 = $ws.Range("A2:A8,A10:A20")
$ws.Range("A2:A8").Font.Name = "Arial"
$ws.Range("A2:A8").Font.Size = 12
$ws.Range("A2:A8").Font.Bold = $true
$ws.Range("A2:A8").HorizontalAlignment = $xlCenter
$ws.Range("A2:A8").VerticalAlignment = $xlCenter

$ws.Range("A10:A20").Font.Name = "Arial"
$ws.Range("A10:A20").Font.Size = 12
$ws.Range("A10:A20").Font.Bold = $true
$ws.Range("A10:A20").HorizontalAlignment = $xlCenter
$ws.Range("A10:A20").VerticalAlignment = $xlCenter

# ---- formatting: column B (siglas), rows 2-8 and 10-20 ----
$ws.Range("B2:B8").Font.Name = "Arial"
$ws.Range("B2:B8").Font.Size = 12
$ws.Range("B2:B8").HorizontalAlignment = $xlRight
$ws.Range("B2:B8").VerticalAlignment = $xlCenter
$ws.Range("B2:B8").IndentLevel = 1

$ws.Range("B10:B20").Font.Name = "Arial"
$ws.Range("B10:B20").Font.Size = 12
$ws.Range("B10:B20").HorizontalAlignment = $xlRight
$ws.Range("B10:B20").VerticalAlignment = $xlCenter
$ws.Range("B10:B20").IndentLevel = 1

# ---- formatting: column C (descriptions), rows 2-8 and 10-20 ----
$ws.Range("C2:C8").Font.Name = "Arial"
$ws.Range("C2:C8").Font.Size = 12
$ws.Range("C2:C8").VerticalAlignment = $xlCenter

$ws.Range("C10:C20").Font.Name = "Arial"
$ws.Range("C10:C20").Font.Size = 12
$ws.Range("C10:C20").VerticalAlignment = $xlCenter

# ---- formatting: spacer rows 9, 21, 24 (column A & B) ----
$spacerA = $ws.Range("A9")
$spacerA.Font.Name = "Arial"
$spacerA.Font.Size = 12
$spacerA.HorizontalAlignment = $xlCenter
$spacerA.VerticalAlignment = $xlCenter

$ws.Range("A21").Font.Name = "Arial"
$ws.Range("A21").Font.Size = 12
$ws.Range("A21").HorizontalAlignment = $xlCenter
$ws.Range("A21").VerticalAlignment = $xlCenter

$ws.Range("A24").Font.Name = "Arial"
$ws.Range("A24").Font.Size = 12
$ws.Range("A24").HorizontalAlignment = $xlCenter
$ws.Range("A24").VerticalAlignment = $xlCenter

$ws.Range("B9").Font.Name = "Arial"
$ws.Range("B9").Font.Size = 12
$ws.Range("B9").HorizontalAlignment = $xlRight
$ws.Range("B9").IndentLevel = 1

$ws.Range("B21").Font.Name = "Arial"
$ws.Range("B21").Font.Size = 12
$ws.Range("B21").HorizontalAlignment = $xlRight
$ws.Range("B21").IndentLevel = 1

$ws.Range("B24").Font.Name = "Arial"
$ws.Range("B24").Font.Size = 12
$ws.Range("B24").HorizontalAlignment = $xlRight
$ws.Range("B24").IndentLevel = 1

$ws.Range("C9").Font.Name = "Arial"
$ws.Range("C9").Font.Size = 12
$ws.Range("C9").HorizontalAlignment = $xlJustify

# ---- formatting: level 4 / level 5 blocks (rows 22-23, 25-26) column A ----
$ws.Range("A22:A23").Font.Name = "Arial"
$ws.Range("A22:A23").Font.Size = 12
$ws.Range("A22:A23").HorizontalAlignment = $xlCenter
$ws.Range("A22:A23").VerticalAlignment = $xlCenter

$ws.Range("A25:A26").Font.Name = "Arial"
$ws.Range("A25:A26").Font.Size = 12
$ws.Range("A25:A26").HorizontalAlignment = $xlCenter
$ws.Range("A25:A26").VerticalAlignment = $xlCenter

# ---- formatting: level 4 / level 5 blocks column B ----
$ws.Range("B22:B23").Font.Name = "Arial"
$ws.Range("B22:B23").Font.Size = 12
$ws.Range("B22:B23").HorizontalAlignment = $xlRight
$ws.Range("B22:B23").IndentLevel = 1

$ws.Range("B25:B26").Font.Name = "Arial"
$ws.Range("B25:B26").Font.Size = 12
$ws.Range("B25:B26").HorizontalAlignment = $xlRight
$ws.Range("B25:B26").IndentLevel = 1

# ---- formatting: level 4 / level 5 blocks column C ----
$ws.Range("C22").Font.Name = "Arial"
$ws.Range("C22").Font.Size = 12

$ws.Range("C23").Font.Name = "Arial"
$ws.Range("C23").Font.Size = 12
$ws.Range("C23").HorizontalAlignment = $xlJustify

$ws.Range("C24").Font.Name = "Arial"
$ws.Range("C24").Font.Size = 12
$ws.Range("C24").HorizontalAlignment = $xlJustify

$ws.Range("C25").Font.Name = "Arial"
$ws.Range("C25").Font.Size = 12
$ws.Range("C25").HorizontalAlignment = $xlJustify

$ws.Range("C26").Font.Name = "Arial"
$ws.Range("C26").Font.Size = 12
$ws.Range("C26").HorizontalAlignment = $xlJustify

# ---- row heights for ht=15.75 rows ----
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(21).RowHeight = 15.75
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 15.75
$ws.Rows.Item(24).RowHeight = 15.75
$ws.Rows.Item(25).RowHeight = 15.75
$ws.Rows.Item(26).RowHeight = 15.75

# ---- merge level-number cells ----
$ws.Range("A2:A8").Merge()
$ws.Range("A10:A20").Merge()
$ws.Range("A22:A23").Merge()
$ws.Range("A25:A26").Merge()

# ---- column widths ----
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# ---- make Plan3 the active sheet / tab ----
$ws.Activate()
$ws.Range("A1:C26").Select()
